$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 934.75
$ws.Range("J17").Value = 969.5
$ws.Range("L17").Value = 2908.5
$ws.Range("N17").Value = -3244.5
$ws.Range("H98").Value = 956.4
$ws.Range("I98").Value = 959.75
$ws.Range("J98").Value = 943
$ws.Range("K98").Value = 959.75
$ws.Range("L98").Value = 943
$ws.Range("M98").Value = 538.25
$ws.Range("N98").Value = -3939
$ws.Range("H122").Value = 956.4
$ws.Range("I122").Value = 959.75
$ws.Range("J122").Value = 943
$ws.Range("K122").Value = 2879.25
$ws.Range("L122").Value = 2829
$ws.Range("M122").Value = -429.25
$ws.Range("N122").Value = -7729
$ws.Range("H132").Value = 2000.8
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 2675.9
$ws.Range("I135").Value = 358.4
$ws.Range("J135").Value = 4993.4
$ws.Range("K135").Value = 3225.6
$ws.Range("L135").Value = 44940.6
$ws.Range("M135").Value = -690.5999999999999
$ws.Range("N135").Value = -50010.6
$ws.Range("H137").Value = 766.3333
$ws.Range("I137").Value = 766.3333
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2298.9999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 251.0001000000002
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 3114.8076
$ws.Range("J138").Value = 3349.5557
$ws.Range("L138").Value = 10048.6671
$ws.Range("N138").Value = -20328.6671

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3250.2144
$ws.Range("I2").Value = 1300.5
$ws.Range("J2").Value = 8124.5
$ws.Range("K2").Value = 1300.5
$ws.Range("L2").Value = 8124.5
$ws.Range("M2").Value = -1187.5
$ws.Range("N2").Value = -8350.5
$ws.Range("H5").Value = 3217.2
$ws.Range("I5").Value = 2521
$ws.Range("K5").Value = 2521
$ws.Range("M5").Value = -2409
$ws.Range("H32").Value = 2894.5173
$ws.Range("I32").Value = 2648.8215
$ws.Range("K32").Value = 2648.8215
$ws.Range("M32").Value = -2361.8215
$ws.Range("H61").Value = 939
$ws.Range("I61").Value = 939
$ws.Range("K61").Value = 939
$ws.Range("M61").Value = -727
$ws.Range("H63").Value = 6000
$ws.Range("I63").Value = 4500
$ws.Range("J63").Value = 6214.2856
$ws.Range("K63").Value = 4500
$ws.Range("L63").Value = 6214.2856
$ws.Range("M63").Value = -3814
$ws.Range("N63").Value = -7586.2856
$ws.Range("H66").Value = 6000
$ws.Range("I66").Value = 4500
$ws.Range("J66").Value = 6214.2856
$ws.Range("K66").Value = 22500
$ws.Range("L66").Value = 31071.428
$ws.Range("M66").Value = -19068
$ws.Range("N66").Value = -37935.428
$ws.Range("H116").Value = 3250.2144
$ws.Range("I116").Value = 1300.5
$ws.Range("J116").Value = 8124.5
$ws.Range("K116").Value = 1300.5
$ws.Range("L116").Value = 8124.5
$ws.Range("M116").Value = 993.5
$ws.Range("N116").Value = -12712.5
$ws.Range("H132").Value = 688.6
$ws.Range("I132").Value = 688.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2065.8
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 464.1999999999998
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 939
$ws.Range("I136").Value = 939
$ws.Range("K136").Value = 2817
$ws.Range("M136").Value = -267

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3250.2144
$ws.Range("I3").Value = 1300.5
$ws.Range("J3").Value = 8124.5
$ws.Range("K3").Value = 1300.5
$ws.Range("L3").Value = 8124.5
$ws.Range("M3").Value = -1186.5
$ws.Range("N3").Value = -8352.5
$ws.Range("H4").Value = 3217.2
$ws.Range("I4").Value = 2521
$ws.Range("K4").Value = 2521
$ws.Range("M4").Value = -2406
$ws.Range("H20").Value = 1644.4286
$ws.Range("I20").Value = 644.2
$ws.Range("K20").Value = 644.2
$ws.Range("M20").Value = -397.2
$ws.Range("H134").Value = 1318.7693
$ws.Range("I134").Value = 1318.7693
$ws.Range("K134").Value = 3956.3079
$ws.Range("M134").Value = -1421.3079

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 855.5
$ws.Range("I16").Value = 855.5
$ws.Range("K16").Value = 855.5
$ws.Range("M16").Value = -568.5
$ws.Range("H113").Value = 855.5
$ws.Range("I113").Value = 855.5
$ws.Range("K113").Value = 855.5
$ws.Range("M113").Value = 1314.5
$ws.Range("H134").Value = 1366.6666
$ws.Range("I134").Value = 1366.6666
$ws.Range("K134").Value = 4099.9998
$ws.Range("M134").Value = -1564.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1787.6
$ws.Range("I122").Value = 1562.7693
$ws.Range("K122").Value = 4688.3079
$ws.Range("M122").Value = -2238.3079
$ws.Range("H132").Value = 4233.1113
$ws.Range("I132").Value = 4233.1113
$ws.Range("K132").Value = 12699.3339
$ws.Range("M132").Value = -10169.3339

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3516
$ws.Range("I122").Value = 1895
$ws.Range("K122").Value = 5685
$ws.Range("M122").Value = -3235
$ws.Range("H132").Value = 5869.769
$ws.Range("I132").Value = 5664.727
$ws.Range("J132").Value = 6997.5
$ws.Range("K132").Value = 16994.181
$ws.Range("L132").Value = 20992.5
$ws.Range("M132").Value = -14464.181
$ws.Range("N132").Value = -26052.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1159.875
$ws.Range("I107").Value = 1076
$ws.Range("K107").Value = 3228
$ws.Range("M107").Value = -1308
$ws.Range("H122").Value = 1033.3334
$ws.Range("I122").Value = 1033.3334
$ws.Range("K122").Value = 3100.0002
$ws.Range("M122").Value = -650.0001999999999
$ws.Range("H126").Value = 1688.909
$ws.Range("I126").Value = 1193.625
$ws.Range("J126").Value = 3009.6667
$ws.Range("K126").Value = 3580.875
$ws.Range("L126").Value = 9029.000100000001
$ws.Range("M126").Value = -1110.875
$ws.Range("N126").Value = -13969.0001
$ws.Range("H132").Value = 1482.5454
$ws.Range("I132").Value = 1482.5454
$ws.Range("K132").Value = 4447.6362
$ws.Range("M132").Value = -1917.6362
